$d = $word.ActiveDocument

# The trailing "_GoBack" bookmark currently sits at the end of the
# paragraph "Add an achievement system of some sorts". The edit relocates
# it into a brand-new, otherwise-empty paragraph appended at the very end
# of the document, after two new brainstorm bullets are added.

# 1) Remove the bookmark from its current location.
$bm = $d.Bookmarks.Item("_GoBack")
$bm.Delete()

# 2) Append the new content at the end of the document:
#      - "Augmented reality extension" as a plain paragraph (default formatting)
#      - an empty paragraph (styled like the others, Times New Roman) that
#        now carries the relocated "_GoBack" bookmark
$endOfDoc = $d.Content.End
$insertionPoint = $d.Range($endOfDoc, $endOfDoc)

$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'
$newParagraphsXml = "<w:p $wNs><w:r><w:t>Augmented reality extension</w:t></w:r></w:p>" `
    + "<w:p $wNs><w:pPr><w:rPr><w:rFonts w:ascii=`"Times New Roman`" w:hAnsi=`"Times New Roman`" w:cs=`"Times New Roman`"/></w:rPr></w:pPr><w:bookmarkStart w:id=`"0`" w:name=`"_GoBack`"/><w:bookmarkEnd w:id=`"0`"/></w:p>"

$insertionPoint.InsertXML($newParagraphsXml)
